# Atualização na porcentagem de conclusão no cronograma
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabela_Tarefas")

# Update completion percentages (cells are formatted as percentage, numFmtId 10 -> 0.00%)
$ws.Range("A2").Value = 0.12
$ws.Range("A3").Value = 0.13
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 1
$ws.Range("A11").Value = 0.12

# Update the selected cell/range shown when the sheet is opened
$ws.Activate()
$ws.Range("C22").Select()
